# Apply the Alvearie FHIR IG deployment update:
#  - Metadata sheet: bump Version, Date, set Publisher value, replace the
#    duplicated "Contact" row with a single "Jurisdiction" row, and delete
#    the now-redundant extra row.
#  - Elements sheet: update the root Extension's Short/Definition text.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# Version: 5.0.0 -> 6.0.0
$metadata.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$metadata.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously empty; now has a publisher name.
$metadata.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact" / "No display for ContactDetail" row
# (rows 10 and 11 were identical); deleting row 11 shifts everything
# below it up by one.
$metadata.Rows.Item(11).Delete()

# The remaining former-"Contact" row (now row 10) becomes a
# "Jurisdiction" / "United States of America" row.
$metadata.Range("A10").Value = "Jurisdiction"
$metadata.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------------

# Root Extension row (row 2): Short & Definition columns (K, L) get the
# resource's actual title/description instead of the generic placeholders.
$elements.Range("K2").Value = "Service Bill Days Count"
$elements.Range("L2").Value = "Number of days between the date of service and the date the claim was received"
